$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")
$ws.Range("Q8").Value = 1
$ws.Range("Q27").Value = 1

$co = $ws.ChartObjects().Add(400, 400, 300, 200)
$ch = $co.Chart
$ch.SetSourceData($ws.Range("F30:T30"))
Write-Host "created test chart"
